# Generate Report for Handoff
# Update the "db3764fe-e09f-4c7a-8a92-c7519ba1e6fa" row with a freshly generated
# handoff timestamp across the Overview summary sheet and each locale sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the db3764fe row (row 7)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-24 12:44:09"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the db3764fe row (row 7)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-24 12:43:57"

# de-de sheet: "Latest Handoff Datetime" column (H) for the db3764fe row (row 7)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-24 12:44:09"
